$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A4 status from N to Y
$ws.Range("A4").Value = "Y"

# Update D4 value from 2 to 3
$ws.Range("D4").Value = 3

# Move the active selection to D4
$ws.Range("D4").Select()

# Set page orientation to portrait
$ws.PageSetup.Orientation = $xlPortrait
